$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 112. Excel shifts the existing rows 112:164
# down to 113:165 and the new row inherits formatting from the row above
# (giving D112 the date style used by the rest of the column).
$ws.Rows.Item(112).Insert()

# Populate the newly inserted row 112 with the new weekly record. All the
# "constant" columns (A, B, C, E, F, G, H, N, O, Q, R) repeat the same
# values used throughout this data block.
$ws.Range("A112").Value = 5
$ws.Range("B112").Value = "Macroferia Regional de Talca"
$ws.Range("C112").Value = "Maule"
$ws.Range("D112").Value = 44510
$ws.Range("E112").Value = 7
$ws.Range("F112").Value = 100112008
$ws.Range("G112").Value = "Coliflor"
$ws.Range("H112").Value = "Sin especificar"
$ws.Range("I112").Value = "Primera"
$ws.Range("J112").Value = 4000
$ws.Range("K112").Value = 600
$ws.Range("L112").Value = 600
$ws.Range("M112").Value = 600
$ws.Range("N112").Value = "`$/unidad"
$ws.Range("O112").Value = "Región del Maule"
$ws.Range("P112").Value = 600
$ws.Range("Q112").Value = 1
$ws.Range("R112").Value = "Hortaliza"
